# Generate Report for Handoff
#
# Adds a new data row (for e31bf546-1abd-403e-9acd-2eea024481ce.md) to the
# "Overview", "zh-cn" and "de-de" report tables, wires up the matching
# hyperlinks, and widens a couple of columns to fit the new content.

$wb = $excel.ActiveWorkbook

$commitSha = "d17c17a5f3ba1c1381af222b684bbe34d84cf3a4"
$newBase   = "e31bf546-1abd-403e-9acd-2eea024481ce"
$newMd     = "$newBase.md"
$ghUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/$commitSha/e2e/$newMd"

function Set-BlankText($range) {
    # Excel drops a cell entirely when you assign "" to it. Writing a lone
    # quote char forces a (empty) literal-text cell, then resetting the
    # style clears the quote-prefix marker Excel adds for that trick.
    $range.Value = "'"
    $range.Style = "Normal"
}

function Style-AsHyperlink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = 15570276
}

function Style-AsDate($range) {
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------
# Sheet "Overview" (sheet1) -> table3 "Overview", columns A-G
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newMd
$wsOverview.Range("B3").Value = "e2e\$newMd"
$wsOverview.Range("C3").Value = ".md"
Set-BlankText $wsOverview.Range("D3")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2017-01-03 04:52:40"
Style-AsDate $wsOverview.Range("G3")

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $ghUrl, "", "", "e2e\$newMd") | Out-Null
Style-AsHyperlink $wsOverview.Range("B3")

$wsOverview.Columns.Item(1).ColumnWidth = 39.084

# ---------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) -> table1 "zh-cn", columns A-R
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A3").Value = $newMd
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = "$newBase.5ed43e904a3e01181c6feb71b05111670583cb61.zh-cn.xlf"
$wsZh.Range("H3").Value = "2017-01-03 04:52:30"
Style-AsDate $wsZh.Range("H3")
Set-BlankText $wsZh.Range("I3")
Set-BlankText $wsZh.Range("J3")
Set-BlankText $wsZh.Range("K3")
$wsZh.Range("L3").Value = "0001-01-01 00:00:00"
Style-AsDate $wsZh.Range("L3")
Set-BlankText $wsZh.Range("M3")
Set-BlankText $wsZh.Range("N3")
$wsZh.Range("O3").Value = "True"
Set-BlankText $wsZh.Range("P3")
$wsZh.Range("Q3").Value = "False"
Set-BlankText $wsZh.Range("R3")

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $ghUrl, "", "", $newMd) | Out-Null
Style-AsHyperlink $wsZh.Range("A3")

$wsZh.Columns.Item(1).ColumnWidth = 39.084
$wsZh.Columns.Item(3).ColumnWidth = 16.25

# ---------------------------------------------------------------------
# Sheet "de-de" (sheet3) -> table2 "de-de", columns A-R
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A3").Value = $newMd
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = "$newBase.5ed43e904a3e01181c6feb71b05111670583cb61.de-de.xlf"
$wsDe.Range("H3").Value = "2017-01-03 04:52:40"
Style-AsDate $wsDe.Range("H3")
Set-BlankText $wsDe.Range("I3")
Set-BlankText $wsDe.Range("J3")
Set-BlankText $wsDe.Range("K3")
$wsDe.Range("L3").Value = "0001-01-01 00:00:00"
Style-AsDate $wsDe.Range("L3")
Set-BlankText $wsDe.Range("M3")
Set-BlankText $wsDe.Range("N3")
$wsDe.Range("O3").Value = "True"
Set-BlankText $wsDe.Range("P3")
$wsDe.Range("Q3").Value = "False"
Set-BlankText $wsDe.Range("R3")

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $ghUrl, "", "", $newMd) | Out-Null
Style-AsHyperlink $wsDe.Range("A3")

$wsDe.Columns.Item(1).ColumnWidth = 39.084
$wsDe.Columns.Item(3).ColumnWidth = 16.25

Write-Host "Report row added for $newMd"
